$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename sheet "Paid" -> "Return"
# ---------------------------------------------------------------------
$paid = $wb.Worksheets.Item("Paid")
$paid.Name = "Return"

# ---------------------------------------------------------------------
# 2. "borrow" sheet updates
# ---------------------------------------------------------------------
$borrow = $wb.Worksheets.Item("borrow")

# Highlight the existing "EarlyPhant" borrow entry (rows 16-17) green,
# marking it returned.
$borrow.Range("A16:AA17").Interior.Color = 5287936

# New plain row (row 20) - a brand-new borrow entry, no special style,
# matching the look of the other un-styled rows (13, 18, 19).
$borrow.Range("A20").Value = 11
$borrow.Range("B20").Value = 7
$borrow.Range("C20").Value = 2016
$borrow.Range("D20").Value = 16
$borrow.Range("E20").Value = "Loem Kimhak"
$borrow.Range("F20").Value = "Kimhak Loem"
$borrow.Range("G20").Value = "Vuthi"
$borrow.Range("K20").Value = 18
$borrow.Range("L20").Value = 7
$borrow.Range("M20").Value = 2016
$borrow.Range("O20").Value = "Rich Dad, Poor Dad - FBSI"

# New merged/styled pair of rows (21-22), mirroring the formatting
# pattern used by rows 14-15 (vertical-center / left+vcenter /
# right+vcenter), recording a return of the "EarlyPhant" book.
$borrow.Range("A21:D22").VerticalAlignment = -4108
$borrow.Range("E21:E22").HorizontalAlignment = -4131
$borrow.Range("E21:E22").VerticalAlignment = -4108
$borrow.Range("G21:G22").HorizontalAlignment = -4131
$borrow.Range("G21:G22").VerticalAlignment = -4108
$borrow.Range("O21:O22").HorizontalAlignment = -4131
$borrow.Range("O21:O22").VerticalAlignment = -4108
$borrow.Range("K21:M22").HorizontalAlignment = -4152
$borrow.Range("K21:M22").VerticalAlignment = -4108
$borrow.Range("H21:J22").HorizontalAlignment = -4108
$borrow.Range("N21:N22").HorizontalAlignment = -4108

$borrow.Range("A21:A22").Merge()
$borrow.Range("B21:B22").Merge()
$borrow.Range("C21:C22").Merge()
$borrow.Range("D21:D22").Merge()
$borrow.Range("E21:E22").Merge()
$borrow.Range("G21:G22").Merge()
$borrow.Range("H21:H22").Merge()
$borrow.Range("I21:I22").Merge()
$borrow.Range("J21:J22").Merge()
$borrow.Range("K21:K22").Merge()
$borrow.Range("L21:L22").Merge()
$borrow.Range("M21:M22").Merge()
$borrow.Range("N21:N22").Merge()
$borrow.Range("O21:O22").Merge()

$borrow.Range("A21").Value = 11
$borrow.Range("B21").Value = 7
$borrow.Range("C21").Value = 2016
$borrow.Range("D21").Value = 17
$borrow.Range("E21").Value = "EarlyPhant"
$borrow.Range("F21").Value = "Marina Moeng"
$borrow.Range("G21").Value = "Vuthi"
$borrow.Range("K21").Value = 18
$borrow.Range("L21").Value = 7
$borrow.Range("M21").Value = 2016
$borrow.Range("O21").Value = "Financial Accounting With Odoo"
$borrow.Range("F22").Value = "Hean Vorthanak"

# ---------------------------------------------------------------------
# 3. "Return" sheet (formerly "Paid") updates
# ---------------------------------------------------------------------
$ret = $wb.Worksheets.Item("Return")
$ret.Range("A13").Value = 11
$ret.Range("B13").Value = 7
$ret.Range("C13").Value = 2016
$ret.Range("D13").Value = 11
$ret.Range("E13").Value = "EarlyPhant"
$ret.Range("F13").Value = "Marina Moeng"
$ret.Range("G13").Value = "Vuthi"
$ret.Range("H13").Value = "Financial Accounting With Odoo"
$ret.Range("E13:H13").Font.Color = $ret.Range("E9").Font.Color
